# "add and more results" - refresh timing/result values and insert a new
# "move_fidelity" summary row ahead of the trailing "Movement times" block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated timing measurements (rows 5-7).
$ws.Range("B5").Value = 0.0006020069122314453
$ws.Range("B6").Value = 0.0006616115570068359
$ws.Range("B7").Value = 0.003393173217773438

# The embedding lists are now emitted as lists-of-lists ("[" / "]") instead
# of lists-of-tuples ("(" / ")") - same 4 rows, text only.
$ws.Range("A8").Value = "[[3, 0], [2, 0], [2, 2], [1, 3], [0, 3], [0, 2], [0, 0], [2, 1], [1, 2], [1, 0], [1, 1], [0, 1]]"
$ws.Range("A49").Value = "[[3, 0], [2, 0], [2, 1], [1, 1], [1, 2], [0, 2], [0, 3], [1, 3], [0, 1], [2, 2], [1, 0], [3, 1]]"
$ws.Range("A108").Value = "[[2, 1], [1, 0], [1, 1], [0, 2], [0, 0], [0, 1], [1, 2], [2, 0], [3, 0], [2, 2], [0, 3], [3, 1]]"
$ws.Range("A154").Value = "[[0, 2], [0, 1], [1, 1], [0, 0], [0, 3], [2, 0], [3, 0], [1, 2], [1, 0], [1, 3], [2, 1], [2, 2]]"

# Insert a new summary row before the existing "Movement times" row (187),
# pushing it and everything below it down by one row.
$ws.Rows.Item(187).Insert()
$ws.Range("A187").Value = "move_fidelity"
$ws.Range("B187").Value = 0.9985968303032804

# The "total time:" row (now row 191, was 190) also reports a fresh value.
$ws.Range("B191").Value = 0.01978707313537598
